# Word COM-interop script implementing:
#   "Updated CAA Record .zone file, and update Cloudflare Guide to
#    include failure reporting URI"
#
# 1. "Configure CAA Records for Cloudflare:" gains a new sentence
#    telling the reader to update the red placeholder text.
# 2. A new CAA "iodef" record (mailto reporting address, highlighted in
#    red) is added as the last bullet under the CAA records list, right
#    before "Go to DNS -> Settings, click ...".
# 3. <w:lastRenderedPageBreak/> markers shift around as a consequence of
#    the extra content pushing the rest of the document down a line.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Find-ParagraphByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Change 1: "Configure CAA Records for Cloudflare:" paragraph gets a new
# run inserted between " for Cloudflare" and ":".
# ---------------------------------------------------------------------
$pCaaIntro = Find-ParagraphByText "*Configure CAA Records for Cloudflare:*"
$rCaaIntro = $pCaaIntro.Range
$rCaaIntro.Collapse(1)
$xmlCaaIntro = "<w:p $wNs>" +
    "<w:pPr><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
    "<w:r><w:t>Configure CAA Records</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> for Cloudflare</w:t></w:r>" +
    "<w:r><w:t>, update the text in red below, this will receive reports on any failed cert issuance</w:t></w:r>" +
    "<w:r><w:t>:</w:t></w:r>" +
    "</w:p>"
$rCaaIntro.InsertXML($xmlCaaIntro)

# ---------------------------------------------------------------------
# Change 2: the last CAA record line (issuewild "ssl.com", right before
# "Go to DNS...") now carries a <w:lastRenderedPageBreak/> on its
# leading "@" run.
# ---------------------------------------------------------------------
$pLastCaa = Find-ParagraphByText '*0 issuewild "ssl.com"*'
$rLastCaa = $pLastCaa.Range
$rLastCaa.Collapse(1)
$xmlLastCaa = "<w:p $wNs>" +
    "<w:pPr><w:numPr><w:ilvl w:val='2'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
    "<w:r><w:lastRenderedPageBreak/><w:t>@</w:t></w:r>" +
    "<w:r><w:tab/><w:t>3600</w:t></w:r>" +
    "<w:r><w:tab/><w:t>CAA</w:t></w:r>" +
    "<w:r><w:tab/><w:t>0 issuewild &quot;ssl.com&quot;</w:t></w:r>" +
    "</w:p>"
$rLastCaa.InsertXML($xmlLastCaa)

# ---------------------------------------------------------------------
# Change 3: insert a brand-new CAA "iodef" record paragraph right before
# "Go to DNS -> Settings, click "Enable DNSSEC"", and drop the page
# break that used to sit on that paragraph (it now renders one bullet
# earlier, see Change 2 above).
# ---------------------------------------------------------------------
$pDns = Find-ParagraphByText "*Go to*DNS*Settings, click*Enable DNSSEC*"
$rDns = $pDns.Range
$rDns.Collapse(1)
$dnsText = $pDns.Range.Text.TrimEnd([char]7)
$xmlNewCaaAndDns = "<w:p $wNs>" +
    "<w:pPr><w:numPr><w:ilvl w:val='2'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
    "<w:r><w:t>@</w:t></w:r>" +
    "<w:r><w:tab/><w:t>3600</w:t></w:r>" +
    "<w:r><w:tab/><w:t>CAA</w:t></w:r>" +
    "<w:r><w:tab/><w:t>0 iodef &quot;</w:t></w:r>" +
    "<w:r><w:rPr><w:color w:val='FF0000'/></w:rPr><w:t>mailto:alert@</w:t></w:r>" +
    "<w:r><w:rPr><w:color w:val='FF0000'/></w:rPr><w:t>domain</w:t></w:r>" +
    "<w:r><w:rPr><w:color w:val='FF0000'/></w:rPr><w:t>.com</w:t></w:r>" +
    "<w:r><w:t>&quot;</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
    "<w:pPr><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
    "<w:r><w:t>Go to" + [char]160 + "DNS -&gt; Settings, click &quot;Enable DNSSEC&quot;</w:t></w:r>" +
    "</w:p>"
$rDns.InsertXML($xmlNewCaaAndDns)

# ---------------------------------------------------------------------
# Change 4: the page break that used to render on "On Caching ->
# Configuration Page" now falls one bullet earlier, on
# "Go to "Protocol Optimization Page"".
# ---------------------------------------------------------------------
$pProtocol = Find-ParagraphByText '*Go to*Protocol Optimization Page*'
$rProtocol = $pProtocol.Range
$rProtocol.Collapse(1)
$xmlProtocol = "<w:p $wNs>" +
    "<w:pPr><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
    "<w:r><w:lastRenderedPageBreak/><w:t>Go to " + [char]8220 + "Protocol Optimization Page" + [char]8221 + "</w:t></w:r>" +
    "</w:p>"
$rProtocol.InsertXML($xmlProtocol)

$pCaching = Find-ParagraphByText "*On Caching -> Configuration Page*"
$rCaching = $pCaching.Range
$rCaching.Collapse(1)
$xmlCaching = "<w:p $wNs>" +
    "<w:pPr><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
    "<w:r><w:t>On Caching -&gt; Configuration Page</w:t></w:r>" +
    "</w:p>"
$rCaching.InsertXML($xmlCaching)

Write-Host "Done."
